# Atualização de bases das ligas, do dia: 28-06-2024 às 19:47
# Swap the stat rows for two pairs of matches (rows 38/39 and 264/265),
# keeping the id column (A) untouched, while all other columns (B..AD)
# swap their contents between the two rows of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol = 30   # column AD

# --- Swap rows 38 and 39 ---
$row1 = 38
$row2 = 39
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $cell1 = $ws.Cells.Item($row1, $c)
    $cell2 = $ws.Cells.Item($row2, $c)
    $v1 = $cell1.Value2
    $v2 = $cell2.Value2
    $cell1.Value2 = $v2
    $cell2.Value2 = $v1
}

# --- Swap rows 264 and 265 ---
$row1 = 264
$row2 = 265
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $cell1 = $ws.Cells.Item($row1, $c)
    $cell2 = $ws.Cells.Item($row2, $c)
    $v1 = $cell1.Value2
    $v2 = $cell2.Value2
    $cell1.Value2 = $v2
    $cell2.Value2 = $v1
}
